# Fix session link in About Me slide.
# Slide 2 ("About Me") has a centered footer-style textbox with the
# tf3604.com session URL; update it from the old "/internals" slug to
# the correct "/temporal" slug.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(10)

$shape.TextFrame.TextRange.Text = "www.tf3604.com/temporal"
